$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E4").Value = "Moose Tracks"
$ws.Range("F4").Value = "PINEAPPLE"
$ws.Range("F4").Select()
